$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Regenerated s_vals data (filtering save games) -- updates B:E and the
# derived "sum" column G for each data row. Column F ("Win") is untouched.

$ws.Range("B2").Value = 0.6606524410359556
$ws.Range("C2").Value = 1.655778082260271
$ws.Range("D2").Value = 0.1494219747398047
$ws.Range("E2").Value = 0.4942365360607697
$ws.Range("G2").Value = 2.960089034096801

$ws.Range("B3").Value = 3.286832544864788
$ws.Range("C3").Value = 10.34677158129881
$ws.Range("D3").Value = 3.537761648806719
$ws.Range("E3").Value = 10.19245300693656
$ws.Range("G3").Value = 27.36381878190688

$ws.Range("B4").Value = 3.286832544864788
$ws.Range("C4").Value = 1.655778082260271
$ws.Range("D4").Value = 3.537761648806719
$ws.Range("E4").Value = 0.4942365360607697
$ws.Range("G4").Value = 8.974608811992548

$ws.Range("B5").Value = 0.6606524410359556
$ws.Range("C5").Value = 1.655778082260271
$ws.Range("D5").Value = 3.537761648806719
$ws.Range("E5").Value = 0.4942365360607697
$ws.Range("G5").Value = 6.348428708163715

$ws.Range("B6").Value = 0.2917716402565462
$ws.Range("C6").Value = 0.306821227259698
$ws.Range("D6").Value = 0.1494219747398047
$ws.Range("E6").Value = 0.4942365360607697
$ws.Range("G6").Value = 1.242251378316819

$ws.Range("B7").Value = 1.455362044514542
$ws.Range("C7").Value = 0.306821227259698
$ws.Range("D7").Value = 0.1494219747398047
$ws.Range("E7").Value = 0.4942365360607697
$ws.Range("G7").Value = 2.405841782574814

$ws.Range("B8").Value = 0.6606524410359556
$ws.Range("C8").Value = 0.306821227259698
$ws.Range("D8").Value = 0.7527432677738641
$ws.Range("E8").Value = 0.4942365360607697
$ws.Range("G8").Value = 2.214453472130288

$ws.Range("B9").Value = 3.286832544864788
$ws.Range("C9").Value = 1.655778082260271
$ws.Range("D9").Value = 3.537761648806719
$ws.Range("E9").Value = 0.4942365360607697
$ws.Range("G9").Value = 8.974608811992548

$ws.Range("B10").Value = 3.286832544864788
$ws.Range("C10").Value = 10.34677158129881
$ws.Range("D10").Value = 0.1494219747398047
$ws.Range("E10").Value = 10.19245300693656
$ws.Range("G10").Value = 23.97547910783996

$ws.Range("B11").Value = 3.286832544864788
$ws.Range("C11").Value = 1.655778082260271
$ws.Range("D11").Value = 3.537761648806719
$ws.Range("E11").Value = 0.4942365360607697
$ws.Range("G11").Value = 8.974608811992548

$ws.Range("B12").Value = 3.286832544864788
$ws.Range("C12").Value = 10.34677158129881
$ws.Range("D12").Value = 3.537761648806719
$ws.Range("E12").Value = 10.19245300693656
$ws.Range("G12").Value = 27.36381878190688
